# Updates the "想去人数" (F column) figures on the "展览" and "全部类型"
# worksheets to reflect newly scraped attendance counts.
# Commit: Update gh-pages to output generated at 456a3b4

$wb = $excel.ActiveWorkbook

# New F-column (column 6) values, keyed by row number, for the "展览" sheet.
$sheet1Updates = @{
    2  = 1048
    3  = 737
    4  = 255
    5  = 26
    6  = 1095
    8  = 1668
    9  = 6117
    11 = 353
    12 = 284
    13 = 87
    14 = 361
    15 = 132
    16 = 5454
    17 = 263
    18 = 1268
    19 = 135
    20 = 112
    22 = 100
    23 = 261
    24 = 96
    26 = 7
    27 = 92
    28 = 1
    29 = 383
    30 = 77
    33 = 42
    34 = 59
    35 = 21
    36 = 61
}

# New F-column (column 6) values, keyed by row number, for the "全部类型" sheet.
$sheet4Updates = @{
    2  = 1048
    3  = 737
    4  = 255
    5  = 26
    6  = 1095
    8  = 1668
    9  = 6118
    11 = 353
    12 = 284
    13 = 87
    14 = 361
    15 = 132
    16 = 5454
    17 = 263
    18 = 1268
    19 = 135
    20 = 112
    22 = 100
    23 = 261
    26 = 7
    27 = 92
    28 = 1
    29 = 383
    30 = 77
    33 = 42
    34 = 59
    35 = 21
    36 = 61
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
